# "Switched to tabbed interface"
# The camera/serial-number table on the sheet was re-sorted: it used to be
# sorted by the "Index" column (A) and is now sorted ascending by the
# "Serial Number" column (C). The header row (row 1) and the extra
# annotation cells in columns D/E/G stay put - only the A:C table body
# (rows 2-27) gets re-ordered, exactly like Excel's own table sort.
# The active selection also moved, from G6 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the table body (A2:C27), leaving the header row alone, by column C
# (Serial Number) ascending - mirrors the Data > Sort operation a user
# would run after clicking into the table and choosing to sort by that
# column.
$dataRange = $ws.Range("A2:C27")
$keyRange = $ws.Range("C2")
$dataRange.Sort($keyRange)

# Reflect the new selection left behind after the sort/tab switch.
$ws.Range("C2").Select()
